# Add 2022-Q4 data
# 1) Create the new "2022-Q4" worksheet by copying the structurally-identical
#    "2022-Q1" sheet (same headers/styles), placed right after "总计".
$wb = $excel.ActiveWorkbook

$srcSheet = $wb.Worksheets.Item("2022-Q1")
$srcSheet.Copy($srcSheet)
$newSheet = $wb.Worksheets.Item("2022-Q1 (2)")
$newSheet.Name = "2022-Q4"

# 2) Fill in the 2022-Q4 fund data (same two funds, new figures).
#    D:G hold numeric-looking text in this workbook (e.g. "21.47"), so force
#    the Text number format first to stop Excel auto-converting them to
#    real numbers (which would also strip significant trailing zeros).
$textCells = $newSheet.Range("D2:G3")
$textCells.NumberFormat = "@"

$newSheet.Range("D2").Value = "21.47"
$newSheet.Range("E2").Value = "92.63"
$newSheet.Range("F2").Value = "3.95"
$newSheet.Range("G2").Value = "0.8481"
$newSheet.Range("H2").Value = 10

$newSheet.Range("D3").Value = "16.17"
$newSheet.Range("E3").Value = "92.63"
$newSheet.Range("F3").Value = "3.95"
$newSheet.Range("G3").Value = "0.6387"
$newSheet.Range("H3").Value = 10

# 3) Update the "总计" (Total) summary sheet: shift the existing rows down
#    by one quarter and insert the new 2022-Q4 figures at the top, plus
#    append the 2021-Q2 row that was previously missing.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("D2").Value = 1.49

$totalSheet.Range("B3").Value = "2022-Q1"
$totalSheet.Range("D3").Value = 2.37

$totalSheet.Range("B4").Value = "2021-Q4"
$totalSheet.Range("D4").Value = 0.84

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 2
$totalSheet.Range("D5").Value = 2.93

# Keep "总计" as the active/selected tab (matches the workbook's original
# bookViews activeTab, which the edit does not otherwise touch).
$totalSheet.Select()
